$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "257.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.08%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.20%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.672"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-10.32%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05881"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.54%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.628"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.78%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8587"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.60%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9494"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-6.44%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.51%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.03961"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "11.50%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07095"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.27%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03178"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.14%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09161"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.75%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001545"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.50%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006028"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.16%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006194"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.28%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.525"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.79%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.203"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.97%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.31%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.97%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.819"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.39%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04224"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001223"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.26%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004303"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.55%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.01%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001936"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "30.47%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03829"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.53%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1103"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.15%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006208"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.90%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002429"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "15.71%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01146"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.37%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005456"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.47%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.01%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06998"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2167"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "9,846.03%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
